$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 91 entirely (the "「100の科学実験」" entry) and shift the rows below it up.
$ws.Rows(91).Delete()
